$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) cells to match the latest scrape.
# NumberFormat is forced to Text ('@') before each write so that values such as
# "0.330" or "11.40" keep their literal text form (matching the original inlineStr
# cells) instead of being auto-coerced into numbers and losing trailing zeros.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '34.838.50'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.838.63'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.77'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.619'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.93%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '39.89'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.330'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0686'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0983'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.105.30'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.40'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.840.37'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.672'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.64'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '34.844.09'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.83'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0786'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '240.35'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.16'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.27'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.38'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.79'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.44'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.123'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.52'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.33%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.95'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.89'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.95%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.04%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +13.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.694'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.07'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '90.72'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.38%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.345.46'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.32%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.94'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.31'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.40'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.76'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.28'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.019.00'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.43'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +21.16%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.73%  '
